$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "95.107.93"
$ws.Range("E2").Value = "  -0.81%  "
$ws.Range("D3").Value = "3.437.15"
$ws.Range("E3").Value = "  +3.90%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'240.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.77%  "
$ws.Range("D6").Value = "'642.19"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.80%  "
$ws.Range("D7").Value = "'1.43"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +7.61%  "
$ws.Range("D8").Value = "'0.404"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.50%  "
$ws.Range("E9").Value = "  +0.10%  "
$ws.Range("D10").Value = "'0.992"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.99%  "
$ws.Range("D11").Value = "3.432.67"
$ws.Range("E11").Value = "  +3.87%  "
$ws.Range("D12").Value = "'41.90"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.67%  "
$ws.Range("D13").Value = "'0.197"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.43%  "
$ws.Range("E14").Value = "  +0.67%  "
$ws.Range("D15").Value = "94.818.81"
$ws.Range("E15").Value = "  -0.96%  "
$ws.Range("D16").Value = "4.085.03"
$ws.Range("E16").Value = "  +3.94%  "
$ws.Range("E17").Value = "  +3.14%  "
$ws.Range("D18").Value = "'8.40"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.19%  "
$ws.Range("D19").Value = "3.442.55"
$ws.Range("E19").Value = "  +3.96%  "
$ws.Range("D20").Value = "'17.98"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +8.03%  "
$ws.Range("D21").Value = "'11.56"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +11.64%  "
$ws.Range("D22").Value = "'0.507"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +7.35%  "
$ws.Range("D23").Value = "'503.11"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.73%  "
$ws.Range("D24").Value = "'3.19"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.75%  "
$ws.Range("D25").Value = "'6.58"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.85%  "
$ws.Range("E26").Value = "  -1.64%  "
$ws.Range("D27").Value = "'91.40"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.54%  "
$ws.Range("E28").Value = "  +2.05%  "
$ws.Range("D29").Value = "3.618.82"
$ws.Range("E29").Value = "  +3.85%  "
$ws.Range("D30").Value = "'11.72"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +9.80%  "
$ws.Range("D31").Value = "'0.998"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.22%  "
$ws.Range("D32").Value = "'2.74"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +12.64%  "
$ws.Range("E33").Value = "  -1.93%  "
$ws.Range("D34").Value = "'0.182"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.27%  "
$ws.Range("D35").Value = "'30.83"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +11.67%  "
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("D37").Value = "'0.565"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.61%  "
$ws.Range("D38").Value = "'7.70"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.05%  "
$ws.Range("D39").Value = "'1.43"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.16%  "
$ws.Range("E40").Value = "  -0.06%  "
$ws.Range("B41").Value = "Bittensor"
$ws.Range("C41").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D41").Value = "'511.24"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.31%  "
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "'0.150"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.00%  "
$ws.Range("E43").Value = "  +11.22%  "
$ws.Range("D44").Value = "'24.04"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.99%  "
$ws.Range("D45").Value = "'1.69"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.82%  "
$ws.Range("D46").Value = "'0.0415"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.31%  "
$ws.Range("D47").Value = "'5.52"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.08%  "
$ws.Range("D48").Value = "'3.54"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.87%  "
$ws.Range("D49").Value = "'2.15"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +10.80%  "
$ws.Range("D50").Value = "'53.63"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.64%  "
$ws.Range("D51").Value = "'3.19"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.71%  "
